$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45202 -> 45203) for every data row (rows 2 through 382).
$ws.Range("C2:C382").Value = 45203
